# Updates odds values on Sheet1 to match the latest FlashScore data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 7 ---
$ws.Range("G7").Value = 1.37
$ws.Range("H7").Value = 4.25
$ws.Range("I7").Value = 8.75
$ws.Range("J7").Value = 1.88
$ws.Range("K7").Value = 2.22
$ws.Range("L7").Value = 7.7
$ws.Range("N7").Value = 7.3
$ws.Range("R7").Value = 1.87
$ws.Range("S7").Value = 1.4
$ws.Range("T7").Value = 2.7
$ws.Range("U7").Value = 2.18
$ws.Range("V7").Value = 1.62
$ws.Range("X7").Value = 5.8
$ws.Range("Z7").Value = 8.25
$ws.Range("AA7").Value = 12
$ws.Range("AC7").Value = 7.3
$ws.Range("AD7").Value = 8.5
$ws.Range("AE7").Value = 23
$ws.Range("AH7").Value = 19
$ws.Range("AI7").Value = 60
$ws.Range("AJ7").Value = 27
$ws.Range("AL7").Value = 120
$ws.Range("AM7").Value = 110
$ws.Range("AN7").Value = 3.05
$ws.Range("AO7").Value = 6.3
$ws.Range("AQ7").Value = 18.5
$ws.Range("AR7").Value = 55
$ws.Range("AT7").Value = 2.7
$ws.Range("AU7").Value = 9
$ws.Range("AX7").Value = 9.25
$ws.Range("AY7").Value = 55
$ws.Range("AZ7").Value = 55
$ws.Range("BA7").Value = 500
$ws.Range("BB7").Value = 500

# --- Row 8 ---
$ws.Range("O8").Value = 1.32
$ws.Range("P8").Value = 2.78

# --- Row 14 ---
$ws.Range("G14").Value = 6.9
$ws.Range("I14").Value = 1.35
$ws.Range("J14").Value = 5.9
$ws.Range("P14").Value = 5.1
$ws.Range("Q14").Value = 1.42
$ws.Range("R14").Value = 2.67
$ws.Range("S14").Value = 1.24
$ws.Range("T14").Value = 3.65
$ws.Range("W14").Value = 25
$ws.Range("AD14").Value = 10.75
$ws.Range("AH14").Value = 10.25
$ws.Range("AL14").Value = 10.25
$ws.Range("AS14").Value = 350
$ws.Range("AT14").Value = 3.65
$ws.Range("AU14").Value = 7.8
$ws.Range("BA14").Value = 14.5
$ws.Range("BC14").Value = 150

# --- Row 17 ---
$ws.Range("Y17").Value = 9
$ws.Range("AQ17").Value = 26
